$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain numeric-looking" string need to be forced to
# Text format first, otherwise Excel COM auto-converts the assigned string into a
# real number (losing the exact decimal text / trailing zeros).
$textCells = @(
    "D5",
    "D6",
    "D10",
    "D11",
    "D14",
    "D17",
    "D20",
    "D22",
    "D23",
    "D25",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D46",
    "D47",
    "D48"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.910.32"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.142.78"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D5").Value = "592.52"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "145.49"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D8").Value = "3.136.00"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").Value = "37.17"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "3.663.31"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "7.31"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "3.143.72"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "63.769.84"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "467.80"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").Value = "81.29"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +5.54%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  +8.04%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  +7.34%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.71"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "27.71"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "2.31"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").Value = "  -5.33%  "
$ws.Range("D40").Value = "51.35"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "456.33"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "9.28"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("D44").Value = "2.922.08"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "39.81"
$ws.Range("E46").Value = "  +10.94%  "
$ws.Range("D47").Value = "0.108"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "129.64"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("E51").Value = "  -1.46%  "
